$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.0498627309343872
$ws.Range("H2").Value = -24.81248665944645
$ws.Range("I2").Value = -3.413942976876101

$ws.Range("G3").Value = 0.1540116713702792
$ws.Range("H3").Value = 30.22823984525034

$ws.Range("G4").Value = -0.3231786052951287
$ws.Range("H4").Value = -17.52147963406326

$ws.Range("G5").Value = -0.3674528162654716
$ws.Range("H5").Value = 7.909992947919488

$ws.Range("G6").Value = 0.1916545661959816
$ws.Range("H6").Value = -2.787034335395721

$ws.Range("G7").Value = 0.3020470976145979
$ws.Range("H7").Value = 45.64749607828841

$ws.Range("G8").Value = 0.1001620517788919
$ws.Range("H8").Value = -1.701701044644506

$ws.Range("G9").Value = 0.1480825482635562
$ws.Range("H9").Value = 17.07975693052233

$ws.Range("G10").Value = 0.06555620062597109
$ws.Range("H10").Value = 6.703357835331801

$ws.Range("G11").Value = 0.04974311481835503
$ws.Range("H11").Value = -0.3745958562658955

$ws.Range("G12").Value = 0.08155946344875906
$ws.Range("H12").Value = -11.89162735155175

$ws.Range("G13").Value = 0.1278482329377386
$ws.Range("H13").Value = 67.76283975217851

$ws.Range("G14").Value = 0.179050365125718
$ws.Range("H14").Value = -20.76072617207298

$ws.Range("G15").Value = 0.2636924799684676
$ws.Range("H15").Value = 7.034313434461096

$ws.Range("G16").Value = 0.1316674038955023
$ws.Range("H16").Value = 15.75679661961901

$ws.Range("G17").Value = 0.1701587347576966
$ws.Range("H17").Value = 13.88834356173248

$ws.Range("G18").Value = -0.02918178016326432
$ws.Range("H18").Value = -225.9824160944987

$ws.Range("G19").Value = 0.04085513655069853
$ws.Range("H19").Value = 68.67101016848125

$ws.Range("G20").Value = 0.1238729350852199
$ws.Range("H20").Value = 45.62808206502136

$ws.Range("G21").Value = 0.08458803709818742
$ws.Range("H21").Value = 29.23286702149963

$ws.Range("G22").Value = 0.1801501288294781
$ws.Range("H22").Value = -5.957115213627667

$ws.Range("G23").Value = 0.2534484925931512
$ws.Range("H23").Value = 17.49647123940774

$ws.Range("G24").Value = -0.01288427524032425
$ws.Range("H24").Value = -238.5637333221544

$ws.Range("G25").Value = 0.004498444500919158
$ws.Range("H25").Value = 119.3420577415382

$ws.Range("G26").Value = 0.2190753399039995
$ws.Range("H26").Value = 6.935751301684147

$ws.Range("G27").Value = 0.2133173888994717
$ws.Range("H27").Value = 10.59330540571808

$ws.Range("G28").Value = 0.02794416365368594
$ws.Range("H28").Value = -58.23848496787092

$ws.Range("G29").Value = 0.09717765027549695
$ws.Range("H29").Value = 3.090315350407567
